{"js": "// 1. \"Contexts (used in same sense as in VC-Data-Model)\" -> \"Context (used in same sense as in VC-Data-Model)\"\nconst body = context.document.body;\n\nconst contextsResults = body.search(\"Contexts (used in same sense as in VC-Data-Model)\", { matchCase: true });\ncontextsResults.load(\"items\");\nawait context.sync();\nif (contextsResults.items.length > 0) {\n  contextsResults.items[0].insertText(\"Context (used in same sense as in VC-Data-Model)\", \"Replace\");\n}\nawait context.sync();\n\n// 2. Insert a new bullet paragraph right before \"first and second item are URIs\"\n//    with the text \"Set the context, which establishes the special terms we will be using\"\nconst uriResults = body.search(\"first and second item are URIs\", { matchCase: true });\nuriResults.load(\"items\");\nawait context.sync();\nif (uriResults.items.length > 0) {\n  const uriParagraph = uriResults.items[0].paragraphs.getFirst();\n  uriParagraph.insertParagraph(\n    \"Set the context, which establishes the special terms we will be using\",\n    \"Before\"\n  );\n}\nawait context.sync();\n\n// 3. \"Diploma Issued as an Open Badge\" -> \"Diploma Issued as an Basic Open Badge\"\nconst diplomaResults = body.search(\"Diploma Issued as an Open Badge\", { matchCase: true });\ndiplomaResults.load(\"items\");\nawait context.sync();\nif (diplomaResults.items.length > 0) {\n  const headingRange = diplomaResults.items[0];\n  const openBadgeResults = headingRange.search(\"Open Badge\", { matchCase: true });\n  openBadgeResults.load(\"items\");\n  await context.sync();\n  if (openBadgeResults.items.length > 0) {\n    openBadgeResults.items[0].insertText(\"Basic \", \"Before\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. \"Contexts (used in same sense as in VC-Data-Model)\" -> \"Context (used in same sense as in VC-Data-Model)\"\n$r1 = $d.Content.Duplicate\nif ($r1.Find.Execute(\"Contexts (used in same sense as in VC-Data-Model)\")) {\n  $r1.Text = \"Context (used in same sense as in VC-Data-Model)\"\n}\n\n# 2. Insert a new bullet paragraph right before \"first and second item are URIs\"\n#    with the text \"Set the context, which establishes the special terms we will be using\"\n$r2 = $d.Content.Duplicate\nif ($r2.Find.Execute(\"first and second item are URIs\")) {\n  $r2.Collapse(1)  # wdCollapseStart\n  $r2.InsertBefore(\"Set the context, which establishes the special terms we will be using`r\")\n}\n\n# 3. \"Diploma Issued as an Open Badge\" -> \"Diploma Issued as an Basic Open Badge\"\n$r3 = $d.Content.Duplicate\nif ($r3.Find.Execute(\"Diploma Issued as an Open Badge\")) {\n  $headingEnd = $r3.End\n  $r4 = $d.Range($r3.Start, $headingEnd)\n  if ($r4.Find.Execute(\"Open Badge\")) {\n    $r4.Collapse(1)  # wdCollapseStart\n    $r4.InsertBefore(\"Basic \")\n  }\n}\n"}
